$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51 (new measurement: Camote, 1a (guarda), from Perú)
$ws.Rows.Item(51).Insert()

$ws.Cells.Item(51,1).Value = 5
$ws.Cells.Item(51,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(51,3).Value = "Maule"
$ws.Cells.Item(51,4).Value = 44477
$ws.Cells.Item(51,5).Value = 7
$ws.Cells.Item(51,6).Value = 100112045
$ws.Cells.Item(51,7).Value = "Zapallo"
$ws.Cells.Item(51,8).Value = "Camote"
$ws.Cells.Item(51,9).Value = "1a (guarda)"
$ws.Cells.Item(51,10).Value = 800
$ws.Cells.Item(51,11).Value = 800
$ws.Cells.Item(51,12).Value = 800
$ws.Cells.Item(51,13).Value = 800
$ws.Cells.Item(51,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(51,15).Value = "Perú"
$ws.Cells.Item(51,16).Value = 800
$ws.Cells.Item(51,17).Value = 1
$ws.Cells.Item(51,18).Value = "Hortaliza"

# Insert a second new row at position 77 (new measurement: Camote, 1a nueva(o))
$ws.Rows.Item(77).Insert()

$ws.Cells.Item(77,1).Value = 5
$ws.Cells.Item(77,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(77,3).Value = "Maule"
$ws.Cells.Item(77,4).Value = 44169
$ws.Cells.Item(77,5).Value = 7
$ws.Cells.Item(77,6).Value = 100112045
$ws.Cells.Item(77,7).Value = "Zapallo"
$ws.Cells.Item(77,8).Value = "Camote"
$ws.Cells.Item(77,9).Value = "1a nueva(o)"
$ws.Cells.Item(77,10).Value = 800
$ws.Cells.Item(77,11).Value = 1000
$ws.Cells.Item(77,12).Value = 1000
$ws.Cells.Item(77,13).Value = 1000
$ws.Cells.Item(77,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(77,15).Value = "Región del Maule"
$ws.Cells.Item(77,16).Value = 1000
$ws.Cells.Item(77,17).Value = 1
$ws.Cells.Item(77,18).Value = "Hortaliza"
